# Commit: "Update gh-pages to output generated at 456a3b4"
#
# This script applies the content refresh to the 上海-漫展信息 workbook:
#  - Sheet "本地生活": the now-expired "次元波板糖×线条小狗MALTESE" pop-up
#    (row 6) is gone from the scrape, so the two rows after it slide up
#    (the serial-number column A keeps its original 5 / 6 values), and the
#    NIJISANJI EN "want-to-go" counter ticks up from 218 to 221.
#  - Sheets "展览", "演出", "本地生活", "全部类型": a handful of "want-to-go"
#    counters (column F) ticked up since the previous scrape.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")

# Row 4 (罗小黑 x HAPPY ZOO) want-to-go counter ticked up.
$wsLocal.Range("F4").Value = 1324

# Drop the expired "次元波板糖×线条小狗MALTESE" row. A full-row delete
# shifts every row below (including column A's serial numbers) up by one,
# so row 7 (PLAVE with animate cafe) becomes the new row 6 and row 8
# (NIJISANJI EN) becomes the new row 7 - content, formatting and cell
# types all carried up intact (no retyping of the date-like text needed).
$wsLocal.Rows.Item(6).Delete()

# The serial-number column (A) must keep its original values (5 and 6)
# rather than the shifted-up ones Delete() just produced.
$wsLocal.Range("A6").Value = 5
$wsLocal.Range("A7").Value = 6

# NIJISANJI EN want-to-go counter ticked up.
$wsLocal.Range("F7").Value = 221

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 4619
$wsExpo.Range("F6").Value = 1791
$wsExpo.Range("F8").Value = 715
$wsExpo.Range("F11").Value = 406
$wsExpo.Range("F12").Value = 1124
$wsExpo.Range("F13").Value = 1566
$wsExpo.Range("F15").Value = 747
$wsExpo.Range("F16").Value = 542
$wsExpo.Range("F19").Value = 150
$wsExpo.Range("F21").Value = 1182
$wsExpo.Range("F23").Value = 2492
$wsExpo.Range("F25").Value = 1526
$wsExpo.Range("F29").Value = 4206

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F18").Value = 278
$wsShow.Range("F22").Value = 243

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All Types)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1324
$wsAll.Range("F5").Value = 221
$wsAll.Range("F8").Value = 4619
$wsAll.Range("F11").Value = 1791
$wsAll.Range("F12").Value = 715
$wsAll.Range("F16").Value = 406
$wsAll.Range("F17").Value = 1124
$wsAll.Range("F18").Value = 1566
$wsAll.Range("F22").Value = 747
$wsAll.Range("F23").Value = 542
$wsAll.Range("F26").Value = 150
$wsAll.Range("F29").Value = 278
$wsAll.Range("F33").Value = 1182
$wsAll.Range("F37").Value = 2492
$wsAll.Range("F38").Value = 243
$wsAll.Range("F43").Value = 1526
$wsAll.Range("F48").Value = 4206
